# Refresh cryptocurrency price/volume(1h) snapshot values (and the
# inserted "OKB" row which shifts rows 9-24 down by one, dropping the
# former last entry) to match the latest scrape, per the GitHub Actions
# commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.146.87"
$ws.Range("E2").Value = "  +0.48%  "
# Row 3
$ws.Range("D3").Value = "1.914.32"
$ws.Range("E3").Value = "  +0.13%  "
# Row 4
$ws.Range("D4").Value = "'1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.47%  "
# Row 5
$ws.Range("D5").Value = "'324.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.11%  "
# Row 6
$ws.Range("E6").Value = "  +0.32%  "
# Row 7
$ws.Range("D7").Value = "'0.4604"
$ws.Range("D7").Style = "Normal"
# Row 8
$ws.Range("D8").Value = "'0.3838"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.80%  "
# Row 9
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "'45.61"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.40%  "
# Row 10
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.07781"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.55%  "
# Row 11
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "'0.9671"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.24%  "
# Row 12
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "'22.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.61%  "
# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.937.33"
$ws.Range("E13").Value = "  +2.55%  "
# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'7.019"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.16%  "
# Row 15
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'5.722"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.77%  "
# Row 16
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "'0.07074"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.28%  "
# Row 17
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "'85.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.18%  "
# Row 18
$ws.Range("B18").Value = "BinanceUSD"
$ws.Range("C18").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D18").Value = "'1.005"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.34%  "
# Row 19
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000009622"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.84%  "
# Row 20
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "'16.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.77%  "
# Row 21
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'1.003"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.27%  "
# Row 22
$ws.Range("B22").Value = "WrappedBTC"
$ws.Range("C22").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").Value = "29.135.12"
$ws.Range("E22").Value = "  +0.39%  "
# Row 23
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "'5.446"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.36%  "
# Row 24
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "'11.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.77%  "
# Row 25
$ws.Range("D25").Value = "'2.087"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.19%  "
# Row 26
$ws.Range("D26").Value = "'157.20"
$ws.Range("D26").Style = "Normal"
# Row 27
$ws.Range("D27").Value = "'19.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.23%  "
# Row 28
$ws.Range("D28").Value = "'5.688"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.64%  "
# Row 29
$ws.Range("D29").Value = "'117.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.37%  "
# Row 30
$ws.Range("D30").Value = "'1.816"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.87%  "
# Row 31
$ws.Range("D31").Value = "'0.09336"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.34%  "
# Row 32
$ws.Range("D32").Value = "'0.8536"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.90%  "
# Row 33
$ws.Range("D33").Value = "'5.103"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.48%  "
# Row 34
$ws.Range("D34").Value = "'1.261"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.71%  "
# Row 35
$ws.Range("D35").Value = "'3.073"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.53%  "
# Row 36
$ws.Range("D36").Value = "'0.05724"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.57%  "
# Row 37
$ws.Range("D37").Value = "'1.159"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.60%  "
# Row 38
$ws.Range("D38").Value = "'0.02065"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.92%  "
# Row 39
$ws.Range("D39").Value = "'7.550"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.16%  "
# Row 40
$ws.Range("D40").Value = "'0.5587"
$ws.Range("D40").Style = "Normal"
# Row 41
$ws.Range("D41").Value = "'0.000003054"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.55%  "
# Row 42
$ws.Range("D42").Value = "'0.1766"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.28%  "
# Row 43
$ws.Range("D43").Value = "'9.170"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.09%  "
# Row 44
$ws.Range("D44").Value = "'2.731"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.86%  "
# Row 45
$ws.Range("D45").Value = "'0.5227"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.63%  "
# Row 46
$ws.Range("D46").Value = "'11.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.02%  "
# Row 47
$ws.Range("D47").Value = "'0.06822"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.47%  "
# Row 48
$ws.Range("D48").Value = "'2.060"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.61%  "
# Row 49
$ws.Range("D49").Value = "'1.793"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.13%  "
# Row 50
$ws.Range("D50").Value = "'110.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.46%  "
# Row 51
$ws.Range("D51").Value = "'0.2977"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.95%  "
